$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the columns that are no longer part of the export (C:I), along with
# their formatting, so the used range shrinks back down to just A:B.
$ws.Range("C1:I5").Clear()

# New header row
$ws.Range("A1").Value = "puntuacion"
$ws.Range("B1").Value = "nombre_completo"

# Keep the score column as literal percentage text (not a numeric % value)
$ws.Range("A2:A11").NumberFormat = "@"

$data = @(
    @("77.19298245614034%", "David Downs"),
    @("82.75862068965517%", "Robert Davis"),
    @("80.0%", "Mary Wells"),
    @("77.41935483870968%", "Bethany Harmon"),
    @("81.9672131147541%", "Brenda Meyers"),
    @("75.0%", "Yvonne Jensen"),
    @("75.86206896551724%", "Joe Ferguson"),
    @("100.0%", "April Gonzalez"),
    @("100.0%", "Colton Collins"),
    @("100.0%", "Joseph Yang")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
